# Auto-generated edit script: updates market-price-derived profit columns (H-N)
# across several worksheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1100.25
$ws.Range("I18").Value = 800.3333
$ws.Range("K18").Value = 800.3333
$ws.Range("M18").Value = -516.3333
$ws.Range("H40").Value = 3561.9656
$ws.Range("I40").Value = 2724.75
$ws.Range("J40").Value = 3880.9048
$ws.Range("K40").Value = 2724.75
$ws.Range("L40").Value = 3880.9048
$ws.Range("M40").Value = -2549.75
$ws.Range("N40").Value = -4230.9048
$ws.Range("H61").Value = 498.7
$ws.Range("I61").Value = 465.33334
$ws.Range("J61").Value = 799
$ws.Range("K61").Value = 1396.00002
$ws.Range("L61").Value = 2397
$ws.Range("M61").Value = -1224.00002
$ws.Range("N61").Value = -2741
$ws.Range("H64").Value = 8997.846
$ws.Range("I64").Value = 7246.75
$ws.Range("K64").Value = 7246.75
$ws.Range("M64").Value = -6998.75
$ws.Range("H67").Value = 8997.846
$ws.Range("I67").Value = 7246.75
$ws.Range("K67").Value = 7246.75
$ws.Range("M67").Value = -6388.75
$ws.Range("H74").Value = 4856.875
$ws.Range("I74").Value = 4477.5
$ws.Range("K74").Value = 4477.5
$ws.Range("M74").Value = -3541.5
$ws.Range("H77").Value = 4856.875
$ws.Range("I77").Value = 4477.5
$ws.Range("K77").Value = 22387.5
$ws.Range("M77").Value = -17707.5
$ws.Range("H135").Value = 3927.0527
$ws.Range("I135").Value = 2539.077
$ws.Range("K135").Value = 22851.693
$ws.Range("M135").Value = -20316.693
$ws.Range("H138").Value = 7655.338
$ws.Range("J138").Value = 7869.422
$ws.Range("L138").Value = 23608.266
$ws.Range("N138").Value = -33888.266

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18150.047
$ws.Range("I32").Value = 17629.834
$ws.Range("K32").Value = 17629.834
$ws.Range("M32").Value = -17342.834
$ws.Range("H61").Value = 3956.1875
$ws.Range("I61").Value = 3735.6428
$ws.Range("K61").Value = 3735.6428
$ws.Range("M61").Value = -3523.6428
$ws.Range("H119").Value = 99964
$ws.Range("J119").Value = 99964
$ws.Range("L119").Value = 99964
$ws.Range("N119").Value = -109640
$ws.Range("H132").Value = 336402.78
$ws.Range("I132").Value = 478432.75
$ws.Range("J132").Value = 4999.4443
$ws.Range("K132").Value = 1435298.25
$ws.Range("L132").Value = 14998.3329
$ws.Range("M132").Value = -1432768.25
$ws.Range("N132").Value = -20058.3329
$ws.Range("H136").Value = 3956.1875
$ws.Range("I136").Value = 3735.6428
$ws.Range("K136").Value = 11206.9284
$ws.Range("M136").Value = -8656.928400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = $null
$ws.Range("H82").Value = 33251
$ws.Range("I82").Value = 25902.2
$ws.Range("K82").Value = 25902.2
$ws.Range("M82").Value = -25519.2
$ws.Range("H85").Value = 33251
$ws.Range("I85").Value = 25902.2
$ws.Range("K85").Value = 25902.2
$ws.Range("M85").Value = -24576.2
$ws.Range("H86").Value = 5106.923
$ws.Range("I86").Value = 4377.4443
$ws.Range("J86").Value = 6748.25
$ws.Range("K86").Value = 4377.4443
$ws.Range("L86").Value = 6748.25
$ws.Range("M86").Value = -3254.4443
$ws.Range("N86").Value = -8994.25
$ws.Range("H89").Value = 5106.923
$ws.Range("I89").Value = 4377.4443
$ws.Range("J89").Value = 6748.25
$ws.Range("K89").Value = 21887.2215
$ws.Range("L89").Value = 33741.25
$ws.Range("M89").Value = -16271.2215
$ws.Range("N89").Value = -44973.25
$ws.Range("H99").Value = 2058.913
$ws.Range("I99").Value = 2049.9375
$ws.Range("K99").Value = 2049.9375
$ws.Range("M99").Value = -551.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 20842.938
$ws.Range("J22").Value = 12672
$ws.Range("L22").Value = 12672
$ws.Range("N22").Value = -13372
$ws.Range("H31").Value = 9788.842000000001
$ws.Range("I31").Value = 4887.2383
$ws.Range("J31").Value = 15843.765
$ws.Range("K31").Value = 4887.2383
$ws.Range("L31").Value = 15843.765
$ws.Range("M31").Value = -4592.2383
$ws.Range("N31").Value = -16433.765
$ws.Range("H34").Value = 9788.842000000001
$ws.Range("I34").Value = 4887.2383
$ws.Range("J34").Value = 15843.765
$ws.Range("K34").Value = 4887.2383
$ws.Range("L34").Value = 15843.765
$ws.Range("M34").Value = -4685.2383
$ws.Range("N34").Value = -16247.765
$ws.Range("H132").Value = 2950.9412
$ws.Range("I132").Value = 2624.4
$ws.Range("K132").Value = 7873.200000000001
$ws.Range("M132").Value = -5343.200000000001
$ws.Range("H134").Value = 6354.5557
$ws.Range("I134").Value = 5824.25
$ws.Range("K134").Value = 17472.75
$ws.Range("M134").Value = -14937.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 106252580
$ws.Range("I4").Value = 67187100
$ws.Range("J4").Value = 399243620
$ws.Range("K4").Value = 201561300
$ws.Range("L4").Value = 1197730860
$ws.Range("M4").Value = -201561188
$ws.Range("N4").Value = -1197731084
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = 2100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 37497.5
$ws.Range("J57").Value = 49995
$ws.Range("L57").Value = 49995
$ws.Range("N57").Value = -51635
$ws.Range("H132").Value = 9363.793
$ws.Range("I132").Value = 9698.261
$ws.Range("J132").Value = 8081.6665
$ws.Range("K132").Value = 29094.783
$ws.Range("L132").Value = 24244.9995
$ws.Range("M132").Value = -26564.783
$ws.Range("N132").Value = -29304.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3285.9167
$ws.Range("I46").Value = 1765.1666
$ws.Range("J46").Value = 4806.6665
$ws.Range("K46").Value = 1765.1666
$ws.Range("L46").Value = 4806.6665
$ws.Range("M46").Value = -1577.1666
$ws.Range("N46").Value = -5182.6665
$ws.Range("H68").Value = 2206.4
$ws.Range("I68").Value = 1933
$ws.Range("K68").Value = 1933
$ws.Range("M68").Value = -1184
$ws.Range("H71").Value = 2206.4
$ws.Range("I71").Value = 1933
$ws.Range("K71").Value = 9665
$ws.Range("M71").Value = -5921
$ws.Range("H122").Value = 5021.5454
$ws.Range("I122").Value = 3747.3333
$ws.Range("J122").Value = 5499.375
$ws.Range("K122").Value = 11241.9999
$ws.Range("L122").Value = 16498.125
$ws.Range("M122").Value = -8791.999899999999
$ws.Range("N122").Value = -21398.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 30006
$ws.Range("I26").Value = 30006
$ws.Range("K26").Value = 30006
$ws.Range("M26").Value = -29713
$ws.Range("H62").Value = 6500
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 6500
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740
$ws.Range("H107").Value = 536.8
$ws.Range("I107").Value = 519.96155
$ws.Range("K107").Value = 1559.88465
$ws.Range("M107").Value = 360.11535
$ws.Range("H122").Value = 33335740
$ws.Range("I122").Value = 58825456
$ws.Range("K122").Value = 176476368
$ws.Range("M122").Value = -176473918
$ws.Range("H124").Value = 65970.22
$ws.Range("J124").Value = 65970.22
$ws.Range("L124").Value = 65970.22
$ws.Range("N124").Value = -75790.22
$ws.Range("H130").Value = 59619.668
$ws.Range("J130").Value = 59619.668
$ws.Range("L130").Value = 59619.668
$ws.Range("N130").Value = -69659.66800000001
$ws.Range("H132").Value = 25596.7
$ws.Range("I132").Value = 39147.723
$ws.Range("K132").Value = 117443.169
$ws.Range("M132").Value = -114913.169
